$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2021 column (R) reusing the formatting of the existing
# 2020 column (Q) so the new cells match the look of the rest of the table.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R4").Value2 = 2021

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R5").Value2 = 72

$excel.CutCopyMode = $false

# Move the selection to the top of the newly added column (this also
# clears the previous "topLeftCell" scroll position of the view).
$ws.Range("R1").Select()
